$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.063527324149735
$ws.Range("D2").Value = 1.064049914034681
$ws.Range("E2").Value = 1.066708009121908
$ws.Range("F2").Value = 1.077131455613958
$ws.Range("I2").Value = 1.056563993208971
$ws.Range("J2").Value = 1.06849224215791
$ws.Range("K2").Value = 1.066766899537356
$ws.Range("L2").Value = 1.069417832507876
$ws.Range("M2").Value = 1.07981356210913
$ws.Range("N2").Value = 1.070009623699548
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.064858264963401
$ws.Range("D3").Value = 1.065119161220602
$ws.Range("E3").Value = 1.068098679100953
$ws.Range("F3").Value = 1.078573836072702
$ws.Range("I3").Value = 1.057076508683191
$ws.Range("J3").Value = 1.069476230069396
$ws.Range("K3").Value = 1.067650639466135
$ws.Range("L3").Value = 1.070622713873603
$ws.Range("M3").Value = 1.081072052812486
$ws.Range("N3").Value = 1.070995008986733
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.065718078392788
$ws.Range("D4").Value = 1.065809762252316
$ws.Range("E4").Value = 1.068995599186635
$ws.Range("F4").Value = 1.079504923101938
$ws.Range("I4").Value = 1.057406084746772
$ws.Range("J4").Value = 1.070111042571849
$ws.Range("K4").Value = 1.068220607169627
$ws.Range("L4").Value = 1.071398888655726
$ws.Range("M4").Value = 1.081883628009349
$ws.Range("N4").Value = 1.071630722995753
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.06607921631991
$ws.Range("D5").Value = 1.066099790442395
$ws.Range("E5").Value = 1.0693719709756
$ws.Range("F5").Value = 1.079895827258425
$ws.Range("I5").Value = 1.057544149278084
$ws.Range("J5").Value = 1.070377469251964
$ws.Range("K5").Value = 1.068459777862467
$ws.Range("L5").Value = 1.071724372127188
$ws.Range("M5").Value = 1.082224163580764
$ws.Range("N5").Value = 1.071897528032314
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.066139833932178
$ws.Range("D6").Value = 1.066148469957565
$ws.Range("E6").Value = 1.069435125069555
$ws.Range("F6").Value = 1.079961431268609
$ws.Range("I6").Value = 1.057567302314606
$ws.Range("J6").Value = 1.070422177302894
$ws.Range("K6").Value = 1.068499909775769
$ws.Range("L6").Value = 1.071778974382398
$ws.Range("M6").Value = 1.082281303087862
$ws.Range("N6").Value = 1.071942299573804
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.065722905215224
$ws.Range("D7").Value = 1.065813638798485
$ws.Range("E7").Value = 1.069000630997725
$ws.Range("F7").Value = 1.079510148433581
$ws.Range("I7").Value = 1.057407931490626
$ws.Range("J7").Value = 1.070114604334536
$ws.Range("K7").Value = 1.068223804718476
$ws.Range("L7").Value = 1.071403240994458
$ws.Range("M7").Value = 1.08188818080864
$ws.Range("N7").Value = 1.071634289816552
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.063977414238139
$ws.Range("D8").Value = 1.064411537432564
$ws.Range("E8").Value = 1.067178606322461
$ws.Range("F8").Value = 1.077619381886152
$ws.Range("I8").Value = 1.056737627613449
$ws.Range("J8").Value = 1.068825180438175
$ws.Range("K8").Value = 1.067065953562493
$ws.Range("L8").Value = 1.06982575042379
$ws.Range("M8").Value = 1.080239449856142
$ws.Range("N8").Value = 1.070343034790359
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06089068073117
$ws.Range("D9").Value = 1.061930905944026
$ws.Range("E9").Value = 1.063945046866503
$ws.Range("F9").Value = 1.074270131925258
$ws.Range("I9").Value = 1.055540591619194
$ws.Range("J9").Value = 1.066538332586749
$ws.Range("K9").Value = 1.065011151077318
$ws.Range("L9").Value = 1.067019080378195
$ws.Range("M9").Value = 1.077312721601917
$ws.Range("N9").Value = 1.068052939352679
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.058825102880752
$ws.Range("D10").Value = 1.060270183693466
$ws.Range("E10").Value = 1.061773327498032
$ws.Range("F10").Value = 1.072024983441084
$ws.Range("I10").Value = 1.05473171007029
$ws.Range("J10").Value = 1.065003566628512
$ws.Range("K10").Value = 1.063631251442606
$ws.Range("L10").Value = 1.065129299650185
$ws.Range("M10").Value = 1.075346621153195
$ws.Range("N10").Value = 1.066515993850741
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.057928757888649
$ws.Range("D11").Value = 1.059549357174866
$ws.Range("E11").Value = 1.060829016016138
$ws.Range("F11").Value = 1.071049762177138
$ws.Range("I11").Value = 1.054378838697881
$ws.Range("J11").Value = 1.064336507957588
$ws.Range("K11").Value = 1.063031299957957
$ws.Range("L11").Value = 1.064306450626646
$ws.Range("M11").Value = 1.074491615185627
$ws.Range("N11").Value = 1.065847987879992
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.057595516935821
$ws.Range("D12").Value = 1.059281345774464
$ws.Range("E12").Value = 1.060477652632511
$ws.Range("F12").Value = 1.070687051319472
$ws.Range("I12").Value = 1.054247369241509
$ws.Range("J12").Value = 1.064088351837853
$ws.Range("K12").Value = 1.062808078509622
$ws.Range("L12").Value = 1.064000111649573
$ws.Range("M12").Value = 1.074173466111916
$ws.Range("N12").Value = 1.065599479350102
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.057667011882922
$ws.Range("D13").Value = 1.059338847180317
$ws.Range("E13").Value = 1.060553048794892
$ws.Range("F13").Value = 1.070764875486798
$ws.Range("I13").Value = 1.054275587924455
$ws.Range("J13").Value = 1.064141599498132
$ws.Range("K13").Value = 1.062855977221469
$ws.Range("L13").Value = 1.064065854081118
$ws.Range("M13").Value = 1.074241735754807
$ws.Range("N13").Value = 1.065652802628167
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05790121822146
$ws.Range("D14").Value = 1.059527208693299
$ws.Range("E14").Value = 1.060799984599871
$ws.Range("F14").Value = 1.071019790041337
$ws.Range("I14").Value = 1.054367979517941
$ws.Range("J14").Value = 1.064316003099006
$ws.Range("K14").Value = 1.063012856038273
$ws.Range("L14").Value = 1.064281142817139
$ws.Range("M14").Value = 1.074465328403496
$ws.Range("N14").Value = 1.065827453902158
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.058045480691257
$ws.Range("D15").Value = 1.059643229245521
$ws.Range("E15").Value = 1.060952049504096
$ws.Range("F15").Value = 1.071176788668452
$ws.Range("I15").Value = 1.054424852270094
$ws.Range("J15").Value = 1.064423408269646
$ws.Range("K15").Value = 1.063109464702204
$ws.Range("L15").Value = 1.064413696730278
$ws.Range("M15").Value = 1.074603016465204
$ws.Range("N15").Value = 1.06593501160046
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.05888454890371
$ws.Range("D16").Value = 1.06031798580918
$ws.Range("E16").Value = 1.06183591427622
$ws.Range("F16").Value = 1.072089640320868
$ws.Range("I16").Value = 1.054755073463696
$ws.Range("J16").Value = 1.065047784040314
$ws.Range("K16").Value = 1.063671016264946
$ws.Range("L16").Value = 1.065183812396235
$ws.Range("M16").Value = 1.075403286861587
$ws.Range("N16").Value = 1.06656027405634
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.059410350497145
$ws.Range("D17").Value = 1.060740777657605
$ws.Range("E17").Value = 1.062389275495818
$ws.Range("F17").Value = 1.072661421853172
$ws.Range("I17").Value = 1.054961508397065
$ws.Range("J17").Value = 1.065438766294596
$ws.Range("K17").Value = 1.064022604181184
$ws.Range("L17").Value = 1.065665656699692
$ws.Range("M17").Value = 1.075904285209295
$ws.Range("N17").Value = 1.066951811550269
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059716855379329
$ws.Range("D18").Value = 1.06098721913711
$ws.Range("E18").Value = 1.062711662119767
$ws.Range("F18").Value = 1.072994638121285
$ws.Range("I18").Value = 1.055081665825936
$ws.Range("J18").Value = 1.065666579317218
$ws.Range("K18").Value = 1.064227443837236
$ws.Range("L18").Value = 1.065946268952531
$ws.Range("M18").Value = 1.076196155515746
$ws.Range("N18").Value = 1.067179948093514
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059821334254625
$ws.Range("D19").Value = 1.061071221346886
$ws.Range("E19").Value = 1.062821523591012
$ws.Range("F19").Value = 1.073108206681665
$ws.Range("I19").Value = 1.055122593672476
$ws.Range("J19").Value = 1.065744217090222
$ws.Range("K19").Value = 1.064297249059193
$ws.Range("L19").Value = 1.066041876281284
$ws.Range("M19").Value = 1.076295616085548
$ws.Range("N19").Value = 1.067257696121061
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.059353956288161
$ws.Range("D20").Value = 1.060695433289559
$ws.Range("E20").Value = 1.06232994442594
$ws.Range("F20").Value = 1.072600105620505
$ws.Range("I20").Value = 1.054939386034614
$ws.Range("J20").Value = 1.06539684250022
$ws.Range("K20").Value = 1.063984906536701
$ws.Range("L20").Value = 1.065614004893716
$ws.Range("M20").Value = 1.075850569473287
$ws.Range("N20").Value = 1.066909828219296
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.057832258594201
$ws.Range("D21").Value = 1.059471748263941
$ws.Range("E21").Value = 1.060727284929435
$ws.Range("F21").Value = 1.070944737128496
$ws.Range("I21").Value = 1.054340783493772
$ws.Range("J21").Value = 1.064264656146633
$ws.Range("K21").Value = 1.062966669469426
$ws.Range("L21").Value = 1.064217764957177
$ws.Range("M21").Value = 1.07439950153171
$ws.Range("N21").Value = 1.065776034031223
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.056873775929841
$ws.Range("D22").Value = 1.058700836788055
$ws.Range("E22").Value = 1.059716128551388
$ws.Range("F22").Value = 1.069901214991442
$ws.Range("I22").Value = 1.053962117736436
$ws.Range("J22").Value = 1.063550599282909
$ws.Range("K22").Value = 1.062324304114241
$ws.Range("L22").Value = 1.063335859957139
$ws.Range("M22").Value = 1.073483902082545
$ws.Range("N22").Value = 1.06506096312485
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.057382052360123
$ws.Range("D23").Value = 1.059109658625871
$ws.Range("E23").Value = 1.060252497445659
$ws.Range("F23").Value = 1.070454667920556
$ws.Range("I23").Value = 1.054163074911904
$ws.Range("J23").Value = 1.063929345491641
$ws.Range("K23").Value = 1.062665040599262
$ws.Range("L23").Value = 1.063803760479056
$ws.Range("M23").Value = 1.073969590675853
$ws.Range("N23").Value = 1.065440247196641
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.05937943897875
$ws.Range("D24").Value = 1.060715922969707
$ws.Range("E24").Value = 1.062356754754401
$ws.Range("F24").Value = 1.072627812694845
$ws.Range("I24").Value = 1.054949382958078
$ws.Range("J24").Value = 1.065415786800353
$ws.Range("K24").Value = 1.064001941207199
$ws.Range("L24").Value = 1.065637345478093
$ws.Range("M24").Value = 1.075874842393715
$ws.Range("N24").Value = 1.066928799422508
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.061690013533042
$ws.Range("D25").Value = 1.062573416229359
$ws.Range("E25").Value = 1.064783777364023
$ws.Range("F25").Value = 1.075138122321007
$ws.Range("I25").Value = 1.055851954236253
$ws.Range("J25").Value = 1.067131312996513
$ws.Range("K25").Value = 1.065544115042375
$ws.Range("L25").Value = 1.067747920701877
$ws.Range("M25").Value = 1.07807194706998
$ws.Range("N25").Value = 1.068646761862641
